# datos.xlsx: add Pruebas / Transfusiones / EfectosAversos sheets,
# renumber Hemocomponentes ids, and refresh the saved selections.

$wb = $excel.ActiveWorkbook

$wsUsuarios = $wb.Worksheets.Item("Usuarios")
$wsHemo     = $wb.Worksheets.Item("Hemocomponentes")
$wsPac      = $wb.Worksheets.Item("Pacientes")

# ---------------------------------------------------------------------
# Hemocomponentes: renumber the id column (10..14 -> 5..9)
# ---------------------------------------------------------------------
$wsHemo.Range("A2").Value = 5
$wsHemo.Range("A3").Value = 6
$wsHemo.Range("A4").Value = 7
$wsHemo.Range("A5").Value = 8
$wsHemo.Range("A6").Value = 9

# ---------------------------------------------------------------------
# New sheet: Pruebas (id, passed)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPruebas = $wb.Worksheets.Add($null, $lastSheet)
$wsPruebas.Name = "Pruebas"

$wsPruebas.Range("A1").Value = "id"
$wsPruebas.Range("B1").Value = "passed"

$wsPruebas.Range("A2").Value = 5
$wsPruebas.Range("B2").Value = "'true"
$wsPruebas.Range("B2").Style = "Normal"

$wsPruebas.Range("A3").Value = 6
$wsPruebas.Range("B3").Value = "'false"
$wsPruebas.Range("B3").Style = "Normal"

$wsPruebas.Range("A4").Value = 7
$wsPruebas.Range("B4").Value = "'true"
$wsPruebas.Range("B4").Style = "Normal"

$wsPruebas.Columns.Item(1).ColumnWidth = 25.5

# ---------------------------------------------------------------------
# New sheet: Transfusiones (hemocomponentId, patientId)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTrans = $wb.Worksheets.Add($null, $lastSheet)
$wsTrans.Name = "Transfusiones"

$wsTrans.Range("A1").Value = "hemocomponentId"
$wsTrans.Range("B1").Value = "patientId"

$wsTrans.Range("A2").Value = 5
$wsTrans.Range("B2").Value = 5

$wsTrans.Range("A3").Value = 6
$wsTrans.Range("B3").Value = 7

$wsTrans.Range("A4").Value = 9
$wsTrans.Range("B4").Value = 6

$wsTrans.Columns.Item(1).ColumnWidth = 19
$wsTrans.Columns.Item(2).ColumnWidth = 27

# ---------------------------------------------------------------------
# New sheet: EfectosAversos (hemocomponentId, patientId, Symptoms)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsEfectos = $wb.Worksheets.Add($null, $lastSheet)
$wsEfectos.Name = "EfectosAversos"

$wsEfectos.Range("A1").Value = "hemocomponentId"
$wsEfectos.Range("B1").Value = "patientId"
$wsEfectos.Range("C1").Value = "Symptoms"

$wsEfectos.Range("A2").Value = 5
$wsEfectos.Range("B2").Value = 5
$wsEfectos.Range("C2").Value = "Fiebre"

$wsEfectos.Range("A3").Value = 9
$wsEfectos.Range("B3").Value = 6
$wsEfectos.Range("C3").Value = "Muchísima fiebre"

$wsEfectos.Columns.Item(3).ColumnWidth = 23

# ---------------------------------------------------------------------
# Restore per-sheet selections
# ---------------------------------------------------------------------
$wsUsuarios.Range("D37").Select()
$wsHemo.Range("A2:A6").Select()
$wsPac.Range("B8").Select()
$wsPruebas.Range("J30").Select()
$wsEfectos.Range("C4").Select()

# Transfusiones is the sheet left active/selected in the saved workbook.
$wsTrans.Range("A1:B4").Select()
$wsTrans.Activate()
